# Update the "dSF" (column F) values for a set of rows in raley_brooks.xlsx
# as part of a re-pull / push of all data and its mean calculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column F (dSF)
$updates = @{
    6  = -3
    8  = -2
    13 = 5
    18 = 2
    19 = 0
    28 = -1
    47 = 0
    48 = -7
    50 = 1
    56 = -8
    59 = 2
    64 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
